$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30; existing rows 30-50 shift down to 31-51.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price entry.
$ws.Cells.Item(30, 1).Value = 10
$ws.Cells.Item(30, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(30, 3).Value = "La Araucanía"
$ws.Cells.Item(30, 4).Value = 44729
$ws.Cells.Item(30, 5).Value = 9
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100107
$ws.Cells.Item(30, 8).Value = "Otros"
$ws.Cells.Item(30, 9).Value = 100107001
$ws.Cells.Item(30, 10).Value = "Caqui"
$ws.Cells.Item(30, 11).Value = "Mankaki"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 30
$ws.Cells.Item(30, 14).Value = 20000
$ws.Cells.Item(30, 15).Value = 20000
$ws.Cells.Item(30, 16).Value = 20000
$ws.Cells.Item(30, 17).Value = '$/bandeja 15 kilos granel'
$ws.Cells.Item(30, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 19).Value = 1333
$ws.Cells.Item(30, 20).Value = 15
